$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency price/volume data scraped on Sat Nov 25 2023.
# All Price (D) and Volume(1h) (E) cells are plain text in the source data
# (e.g. "37.879.61" and padded percentages like "  +0.08%  "), so values
# that Excel would otherwise auto-parse as numbers are forced to stay text
# by temporarily applying a Text number format, then restoring the default
# "Normal" cell style so no visible formatting change is introduced.

$ws.Range("D2").Value = "37.879.61"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "2.083.66"
$ws.Range("E3").Value = "  -0.11%  "
$ws.Range("E4").Value = "  +0.09%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "233.31"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("E6").Value = "  +0.31%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "59.29"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +3.41%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +2.04%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.0787"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +1.12%  "
$ws.Range("E11").Value = "  +1.38%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "14.73"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +2.49%  "
$ws.Range("E13").Value = "  +0.83%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.776"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +1.89%  "
$ws.Range("E15").Value = "  +2.29%  "
$ws.Range("D16").Value = "2.082.08"
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("D17").Value = "37.814.49"
$ws.Range("E17").Value = "  +0.20%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "6.12"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -0.27%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "71.65"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +1.08%  "
$ws.Range("D20").Value = "0.0₃0847"
$ws.Range("E20").Value = "  +3.17%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "228.08"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("E23").Value = "  -0.47%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "2.39"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.72%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "172.00"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.00%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "9.21"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +3.30%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.137"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -1.28%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "1.41"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -2.06%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "19.48"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("E30").Value = "  +1.81%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "4.72"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +2.48%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "4.71"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +2.86%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.0631"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +1.03%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "2.50"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +0.19%  "
$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.82"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -0.71%  "
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "3.41"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +0.68%  "
$ws.Range("E37").Value = "  -0.06%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "5.41"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +0.13%  "
$ws.Range("E39").Value = "  -1.12%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "99.02"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +2.07%  "
$ws.Range("E41").Value = "  +2.38%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "2.89"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -1.25%  "
$ws.Range("E43").Value = "  +8.26%  "
$ws.Range("D44").Value = "1.444.70"
$ws.Range("E44").Value = "  -0.50%  "
$ws.Range("E45").Value = "  -0.52%  "
$ws.Range("E46").Value = "  +2.87%  "
$ws.Range("E47").Value = "  +0.88%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "7.37"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +0.06%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "2.99"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -0.09%  "
$ws.Range("D50").Value = "2.276.97"
$ws.Range("E50").Value = "  -0.02%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "46.79"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.75%  "
